# Rename "Device" sheet to "Apparatus" and update the "Device" wording
# used in its cell text to "Apparatus", matching the commit:
#   Change "Device" to "Apparatus" in excel form, simulink, function name

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# Update cell text that referred to "device"/"Device" wording.
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."
$ws.Range("A8").Value = "Apparatus type with default values:"
$ws.Range("B9").Value = "Apparatus type"
$ws.Range("B56").Value = "Apparatus type"
$ws.Range("C56").Value = "Apparatus parameters"

# Rename the sheet itself last (name lookups above use the old name).
$ws.Name = "Apparatus"

# Reset the view to the top-left cell with no special selection, matching
# the saved state in the target workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
